$wb = $excel.ActiveWorkbook

# --- 1. Create "week 3" as a copy of "week 2", placed right after it ---
$week2 = $wb.Worksheets.Item("week 2")
$week2.Copy($null, $week2) | Out-Null
$ws3 = $wb.Worksheets.Item("week 2 (2)")
$ws3.Name = "week 3"

# Clear last week's log entries (rows 8-14) but keep the formatting/formulas
$ws3.Range("C8:F14").ClearContents()

# Fill in this week's first log entry (row 7)
$ws3.Range("A7").Value = "Woensdag"
$ws3.Range("B7").Value = 41661
$ws3.Range("B7").NumberFormat = "mm-dd-yy"
$ws3.Range("B7").HorizontalAlignment = -4108
$ws3.Range("C7").Value = 0.61319444444444449
$ws3.Range("D7").Value = 0.61805555555555558
$ws3.Range("F7").Value = "beetle.cs en ibeetlestate.cs aangemaakt"

# Resize rows to fit their (now shorter) content again
$ws3.Range("A7:A14").EntireRow.AutoFit()

$ws3.Range("F7").Select() | Out-Null

# --- 2. Add a "week 3" line to the "Totaal" overview sheet ---
$total = $wb.Worksheets.Item("Totaal")
$total.Range("A11").EntireRow.Insert()
$total.Range("A11").Value = 3
$total.Range("B11").Formula = "='week 3'!G18"
$total.Range("B12").Select() | Out-Null
